# Applies the Sheets update scraped by the scheduled runner:
# updates currentAveragePrice / LevePrice / LeveProfit figures for the
# affected Leve rows across the ALC, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 728117.2
$ws.Range("I2").Value = 1212396.4
$ws.Range("K2").Value = 1212396.4
$ws.Range("M2").Value = -1212283.4

$ws.Range("H62").Value = 10166.071
$ws.Range("I62").Value = 12281.9
$ws.Range("K62").Value = 12281.9
$ws.Range("M62").Value = -11657.9

$ws.Range("H65").Value = 10166.071
$ws.Range("I65").Value = 12281.9
$ws.Range("K65").Value = 61409.5
$ws.Range("M65").Value = -58289.5

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H138").Value = 2303.4146
$ws.Range("J138").Value = 2987.375
$ws.Range("L138").Value = 8962.125
$ws.Range("N138").Value = -19242.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5266145
$ws.Range("I20").Value = 7693343.5
$ws.Range("J20").Value = 7215.5
$ws.Range("K20").Value = 7693343.5
$ws.Range("L20").Value = 7215.5
$ws.Range("M20").Value = -7693096.5
$ws.Range("N20").Value = -7709.5

$ws.Range("H55").Value = 74999
$ws.Range("J55").Value = 74999
$ws.Range("L55").Value = 74999
$ws.Range("N55").Value = -75545

$ws.Range("H105").Value = 2657.4443
$ws.Range("I105").Value = 2671.923
$ws.Range("J105").Value = 2619.8
$ws.Range("K105").Value = 2671.923
$ws.Range("L105").Value = 2619.8
$ws.Range("M105").Value = -924.9229999999998
$ws.Range("N105").Value = -6113.8

$ws.Range("H134").Value = 2645.6765
$ws.Range("I134").Value = 2460.9167
$ws.Range("K134").Value = 7382.750100000001
$ws.Range("M134").Value = -4847.750100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 106.09091
$ws.Range("I7").Value = 142.8
$ws.Range("J7").Value = 75.5
$ws.Range("K7").Value = 142.8
$ws.Range("L7").Value = 75.5
$ws.Range("M7").Value = -29.80000000000001
$ws.Range("N7").Value = -301.5

$ws.Range("H31").Value = 5771.394
$ws.Range("I31").Value = 4716.5
$ws.Range("K31").Value = 4716.5
$ws.Range("M31").Value = -4421.5

$ws.Range("H34").Value = 5771.394
$ws.Range("I34").Value = 4716.5
$ws.Range("K34").Value = 4716.5
$ws.Range("M34").Value = -4514.5

$ws.Range("H50").Value = 46974
$ws.Range("J50").Value = 46974
$ws.Range("L50").Value = 46974
$ws.Range("N50").Value = -48224

$ws.Range("H122").Value = 2439.111
$ws.Range("I122").Value = 2495.5
$ws.Range("J122").Value = 1988
$ws.Range("K122").Value = 7486.5
$ws.Range("L122").Value = 5964
$ws.Range("M122").Value = -5036.5
$ws.Range("N122").Value = -10864

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 2239.6
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 2239.6
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 6718.799999999999
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -11960.8

$ws.Range("H113").Value = 810.6667
$ws.Range("I113").Value = 810.8570999999999
$ws.Range("K113").Value = 2432.5713
$ws.Range("M113").Value = -262.5712999999996

$ws.Range("H115").Value = 28
$ws.Range("I115").Value = 28
$ws.Range("K115").Value = 84
$ws.Range("M115").Value = 1091

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 76429.664
$ws.Range("J42").Value = 89644.5
$ws.Range("L42").Value = 89644.5
$ws.Range("N42").Value = -90614.5

$ws.Range("H115").Value = 76429.664
$ws.Range("J115").Value = 89644.5
$ws.Range("L115").Value = 89644.5
$ws.Range("N115").Value = -91994.5

$ws.Range("H122").Value = 2814.6155
$ws.Range("I122").Value = 2190
$ws.Range("K122").Value = 6570
$ws.Range("M122").Value = -4120

$ws.Range("H132").Value = 44484.58
$ws.Range("I132").Value = 47608.082
$ws.Range("K132").Value = 142824.246
$ws.Range("M132").Value = -140294.246

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5631.4443
$ws.Range("J40").Value = 8000
$ws.Range("L40").Value = 8000
$ws.Range("N40").Value = -8272

$ws.Range("H46").Value = 6957.4165
$ws.Range("I46").Value = 37998.668
$ws.Range("J46").Value = 2522.9524
$ws.Range("K46").Value = 37998.668
$ws.Range("L46").Value = 2522.9524
$ws.Range("M46").Value = -37810.668
$ws.Range("N46").Value = -2898.9524

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H74").Value = 79877.336
$ws.Range("J74").Value = 79877.336
$ws.Range("L74").Value = 79877.336
$ws.Range("N74").Value = -81873.336

$ws.Range("H77").Value = 79877.336
$ws.Range("J77").Value = 79877.336
$ws.Range("L77").Value = 239632.008
$ws.Range("N77").Value = -249616.008

$ws.Range("H136").Value = 4327.28
$ws.Range("J136").Value = 8080
$ws.Range("L136").Value = 24240
$ws.Range("N136").Value = -29340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 19000
$ws.Range("J68").Value = 19000
$ws.Range("L68").Value = 19000
$ws.Range("N68").Value = -20622

$ws.Range("H71").Value = 19000
$ws.Range("J71").Value = 19000
$ws.Range("L71").Value = 57000
$ws.Range("N71").Value = -65112

$ws.Range("H122").Value = 16212.692
$ws.Range("I122").Value = 18177.6
$ws.Range("J122").Value = 9663
$ws.Range("K122").Value = 54532.8
$ws.Range("L122").Value = 28989
$ws.Range("M122").Value = -52082.8
$ws.Range("N122").Value = -33889

$ws.Range("H132").Value = 32034.516
$ws.Range("I132").Value = 32034.516
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 96103.548
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -93573.548
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 2104.3242
$ws.Range("I136").Value = 1678.258
$ws.Range("J136").Value = 4305.6665
$ws.Range("K136").Value = 4305.6665
$ws.Range("L136").Value = 12916.9995
$ws.Range("M136").Value = -2484.774
$ws.Range("N136").Value = -18016.9995
